# Auto-generated edit script: refresh market price data across all Leve profit sheets
# Mirrors a scheduled market-data refresh run (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 58.75
$ws.Range("I12").Value = 58.75
$ws.Range("K12").Value = 58.75
$ws.Range("M12").Value = 111.25
$ws.Range("H113").Value = 6333.3335
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 7500
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 7500
$ws.Range("M113").Value = -746
$ws.Range("N113").Value = -14008
$ws.Range("H137").Value = 7371.0835
$ws.Range("I137").Value = 5700
$ws.Range("J137").Value = 9042.166999999999
$ws.Range("K137").Value = 17100
$ws.Range("L137").Value = 27126.501
$ws.Range("M137").Value = -14550
$ws.Range("N137").Value = -32226.501

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H22").Value = 1999.5
$ws.Range("I22").Value = 1999.5
$ws.Range("K22").Value = 1999.5
$ws.Range("M22").Value = -1700.5
$ws.Range("H32").Value = 9100.666999999999
$ws.Range("I32").Value = 2901
$ws.Range("K32").Value = 2901
$ws.Range("M32").Value = -2614
$ws.Range("H45").Value = 6609.875
$ws.Range("I45").Value = 3137.182
$ws.Range("J45").Value = 14249.8
$ws.Range("K45").Value = 3137.182
$ws.Range("L45").Value = 14249.8
$ws.Range("M45").Value = -2760.182
$ws.Range("N45").Value = -15003.8
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 15000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 15000
$ws.Range("N46").Value = -15638
$ws.Range("H102").Value = 2474.5
$ws.Range("I102").Value = 2474.5
$ws.Range("K102").Value = 2474.5
$ws.Range("M102").Value = -852.5
$ws.Range("H132").Value = 8492.25
$ws.Range("I132").Value = 2787.6
$ws.Range("J132").Value = 18000
$ws.Range("K132").Value = 8362.799999999999
$ws.Range("L132").Value = 54000
$ws.Range("M132").Value = -5832.799999999999
$ws.Range("N132").Value = -59060
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1357.2
$ws.Range("I11").Value = 264.33334
$ws.Range("J11").Value = 2996.5
$ws.Range("K11").Value = 264.33334
$ws.Range("L11").Value = 2996.5
$ws.Range("M11").Value = -124.33334
$ws.Range("N11").Value = -3276.5
$ws.Range("H94").Value = 1271.3636
$ws.Range("I94").Value = 898.125
$ws.Range("K94").Value = 898.125
$ws.Range("M94").Value = -447.125
$ws.Range("H107").Value = 2999.6667
$ws.Range("I107").Value = 1999.5
$ws.Range("K107").Value = 1999.5
$ws.Range("M107").Value = -79.5
$ws.Range("H134").Value = 5910.25
$ws.Range("I134").Value = 2592.3
$ws.Range("K134").Value = 7776.900000000001
$ws.Range("M134").Value = -5241.900000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 817.6667
$ws.Range("I22").Value = 1102
$ws.Range("J22").Value = 533.3333
$ws.Range("K22").Value = 1102
$ws.Range("L22").Value = 533.3333
$ws.Range("M22").Value = -752
$ws.Range("N22").Value = -1233.3333
$ws.Range("H41").Value = 4000
$ws.Range("I41").Value = 4000
$ws.Range("K41").Value = 4000
$ws.Range("M41").Value = -3572
$ws.Range("H92").Value = 44999.5
$ws.Range("J92").Value = 44999.5
$ws.Range("L92").Value = 44999.5
$ws.Range("N92").Value = -49991.5
$ws.Range("H99").Value = 5846.4287
$ws.Range("I99").Value = 6182.2
$ws.Range("J99").Value = 5007
$ws.Range("K99").Value = 6182.2
$ws.Range("L99").Value = 5007
$ws.Range("M99").Value = -4684.2
$ws.Range("N99").Value = -8003
$ws.Range("H105").Value = 1519.6
$ws.Range("I105").Value = 1274.5
$ws.Range("K105").Value = 1274.5
$ws.Range("M105").Value = 472.5
$ws.Range("H126").Value = 5846.4287
$ws.Range("I126").Value = 6182.2
$ws.Range("J126").Value = 5007
$ws.Range("K126").Value = 18546.6
$ws.Range("L126").Value = 15021
$ws.Range("M126").Value = -16076.6
$ws.Range("N126").Value = -19961
$ws.Range("H134").Value = 8432.691999999999
$ws.Range("I134").Value = 3857.5
$ws.Range("K134").Value = 11572.5
$ws.Range("M134").Value = -9037.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 274.5
$ws.Range("I41").Value = 274.5
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 823.5
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -485.5
$ws.Range("H55").Value = 1875.5
$ws.Range("J55").Value = 1875.5
$ws.Range("L55").Value = 5626.5
$ws.Range("N55").Value = -5980.5
$ws.Range("H92").Value = 497
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H108").Value = 5258
$ws.Range("J108").Value = 9989
$ws.Range("L108").Value = 29967
$ws.Range("N108").Value = -35727
$ws.Range("N41").ClearContents()
$ws.Range("N92").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 50000000
$ws.Range("I14").Value = 50000000
$ws.Range("K14").Value = 50000000
$ws.Range("M14").Value = -49999832
$ws.Range("H36").Value = 1100
$ws.Range("J36").Value = 1100
$ws.Range("L36").Value = 1100
$ws.Range("N36").Value = -2070
$ws.Range("H43").Value = 7084.615
$ws.Range("I43").Value = 1000
$ws.Range("J43").Value = 14183.333
$ws.Range("K43").Value = 1000
$ws.Range("L43").Value = 14183.333
$ws.Range("M43").Value = -849
$ws.Range("N43").Value = -14485.333
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 1500
$ws.Range("K80").Value = 1500
$ws.Range("M80").Value = -502
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 1500
$ws.Range("K83").Value = 7500
$ws.Range("M83").Value = -2508
$ws.Range("H102").Value = 1183.5
$ws.Range("I102").Value = 1183.5
$ws.Range("K102").Value = 1183.5
$ws.Range("M102").Value = 438.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4900
$ws.Range("J16").Value = 4900
$ws.Range("L16").Value = 4900
$ws.Range("N16").Value = -5240
$ws.Range("H22").Value = 1340.4
$ws.Range("I22").Value = 1666.6666
$ws.Range("K22").Value = 1666.6666
$ws.Range("M22").Value = -1371.6666
$ws.Range("H27").Value = 1340.4
$ws.Range("I27").Value = 1666.6666
$ws.Range("K27").Value = 1666.6666
$ws.Range("M27").Value = -1559.6666
$ws.Range("H46").Value = 943
$ws.Range("I46").Value = 1115.6666
$ws.Range("J46").Value = 425
$ws.Range("K46").Value = 1115.6666
$ws.Range("L46").Value = 425
$ws.Range("M46").Value = -927.6666
$ws.Range("N46").Value = -801
$ws.Range("H55").Value = 1394.6
$ws.Range("I55").Value = 1394.6
$ws.Range("K55").Value = 1394.6
$ws.Range("M55").Value = -1221.6
$ws.Range("H82").Value = 2030.1
$ws.Range("J82").Value = 2050.3333
$ws.Range("L82").Value = 2050.3333
$ws.Range("N82").Value = -2772.3333
$ws.Range("H85").Value = 2030.1
$ws.Range("J85").Value = 2050.3333
$ws.Range("L85").Value = 2050.3333
$ws.Range("N85").Value = -4546.3333
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 500500
$ws.Range("I29").Value = 500500
$ws.Range("K29").Value = 500500
$ws.Range("M29").Value = -500210
$ws.Range("H122").Value = 998.25
$ws.Range("I122").Value = 997.6667
$ws.Range("K122").Value = 2993.0001
$ws.Range("M122").Value = -543.0001000000002

Write-Output "Updated 207 cells across 8 sheets"